$d = $word.ActiveDocument

$replacements = @(
    @("66×88=", "68×73="),
    @("40×21=", "54×71="),
    @("99×24=", "29×88="),
    @("37×37=", "15×67="),
    @("12×62=", "29×97="),
    @("88×78=", "54×11="),
    @("50×68=", "55×51="),
    @("73×47=", "67×55="),
    @("94×49=", "58×78="),
    @("37×94=", "18×99="),
    @("15×15=", "55×83="),
    @("33×56=", "27×43="),
    @("82×83=", "64×52="),
    @("48×22=", "96×36="),
    @("43×89=", "17×47="),
    @("68×86=", "34×60="),
    @("81×52=", "15×13="),
    @("47×68=", "56×91="),
    @("39×71=", "77×46="),
    @("69×90=", "50×34="),
    @("87×27=", "66×42="),
    @("54×83=", "89×53="),
    @("58×94=", "29×67="),
    @("41×42=", "83×75="),
    @("66×40=", "92×95=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
